$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update simple scalar values on the Metadata sheet ---
$ws1.Range("B3").Value = "0.1.7"
$ws1.Range("B6").Value = "draft"
$ws1.Range("B8").Value = "2024-11-22T12:33:30-06:00"
$ws1.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws1.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- Insert a new "Jurisdiction" row after the Contact rows (new row 12) ---
# First extend the existing formatting down to the new row 16 by copying the
# format of row 15 (avoids creating brand-new style entries via Rows.Insert()).
$ws1.Range("A15:B15").Copy()
$ws1.Range("A16:B16").PasteSpecial(-4122)

# Shift the Description/Purpose/Copyright/Immutable rows down by one to make
# room for the new Jurisdiction row.
$ws1.Range("A16").Value = $ws1.Range("A15").Value2
$ws1.Range("B16").Value = $ws1.Range("B15").Value2
$ws1.Range("A15").Value = $ws1.Range("A14").Value2
$ws1.Range("B15").Value = $ws1.Range("B14").Value2
$ws1.Range("A14").Value = $ws1.Range("A13").Value2
$ws1.Range("B14").Value = $ws1.Range("B13").Value2
$ws1.Range("A13").Value = $ws1.Range("A12").Value2
$ws1.Range("B13").Value = $ws1.Range("B12").Value2

# Fill in the new Jurisdiction row.
$ws1.Range("A12").Value = "Jurisdiction"
$ws1.Range("B12").Value = ""
